$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.683.40"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "2.924.18"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'548.12"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'129.93"
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "2.916.47"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'32.61"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "3.406.09"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "'6.85"
$ws.Range("E17").Value = "  +5.39%  "
$ws.Range("D18").Value = "2.918.89"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "57.658.15"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "'415.48"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("D23").Value = "'13.29"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("D24").Value = "'6.94"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "'79.39"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'2.46"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "'7.33"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("D31").Value = "'25.15"
$ws.Range("E31").Value = "  +0.35%  "
$ws.Range("D32").Value = "'5.93"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").Value = "'0.0962"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("E36").Value = "  +4.01%  "
$ws.Range("D37").Value = "'48.13"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0686"
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").Value = "'8.74"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  +6.81%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "2.701.19"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").Value = "'370.90"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D46").Value = "'123.80"
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'1.93"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'22.68"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -0.68%  "
